$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8510355029585798
$ws.Range("C2").Value = 0.1012245750502546

$ws.Range("B3").Value = 0.7795857988165681
$ws.Range("C3").Value = 0.09496585691618256

$ws.Range("B4").Value = 0.8126725838264299
$ws.Range("C4").Value = 0.1780890680725701

$ws.Range("B5").Value = 0.6672090729783038
$ws.Range("C5").Value = 0.1677474570688943

$ws.Range("B6").Value = 0.8407914201183431
$ws.Range("C6").Value = 0.06079329372474069
